# Update cryptos list figures (Price / Volume(1h)) per the scheduled data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @{ D = <new Price text>; E = <new Volume(1h) text> }
# D is omitted for rows where only the Volume(1h) figure changed.
$updates = @{
    2  = @{ D = "26.851.62"; E = "  -1.87%  " }
    3  = @{ D = "1.812.56";  E = "  -0.51%  " }
    4  = @{ D = "1.002";     E = "  +0.26%  " }
    5  = @{ D = "310.27";    E = "  -1.12%  " }
    6  = @{ E = "  +0.17%  " }
    7  = @{ D = "0.4631";    E = "  -0.61%  " }
    8  = @{ E = "  -2.10%  " }
    9  = @{ D = "0.07346";   E = "  -1.24%  " }
    10 = @{ D = "0.8695";    E = "  -0.30%  " }
    11 = @{ D = "20.39";     E = "  -1.41%  " }
    12 = @{ D = "1.904.14";  E = "  +4.46%  " }
    13 = @{ E = "  -1.30%  " }
    14 = @{ D = "0.07075";   E = "  -0.40%  " }
    15 = @{ E = "  -2.59%  " }
    16 = @{ D = "91.65";     E = "  -0.63%  " }
    17 = @{ D = "1.001";     E = "  +0.16%  " }
    18 = @{ D = "0.000008699"; E = "  -0.86%  " }
    19 = @{ D = "1.001";     E = "  +0.12%  " }
    20 = @{ D = "14.66";     E = "  -1.98%  " }
    21 = @{ D = "26.892.29"; E = "  -1.73%  " }
    22 = @{ D = "5.331";     E = "  +0.39%  " }
    23 = @{ D = "10.55";     E = "  -3.39%  " }
    24 = @{ D = "2.065.38";  E = "  +0.63%  " }
    25 = @{ E = "  -2.25%  " }
    26 = @{ D = "151.87";    E = "  +0.22%  " }
    27 = @{ D = "18.35";     E = "  -1.47%  " }
    28 = @{ D = "2.120";     E = "  -6.22%  " }
    29 = @{ E = "  -0.12%  " }
    30 = @{ D = "115.43";    E = "  -1.43%  " }
    31 = @{ D = "0.08890";   E = "  -0.08%  " }
    32 = @{ D = "0.7559";    E = "  -3.54%  " }
    33 = @{ D = "2.932";     E = "  +0.39%  " }
    34 = @{ D = "1.150";     E = "  -2.95%  " }
    35 = @{ D = "4.457";     E = "  -1.72%  " }
    36 = @{ E = "  +0.25%  " }
    37 = @{ D = "1.093";     E = "  -0.46%  " }
    38 = @{ D = "0.01954";   E = "  -1.04%  " }
    39 = @{ D = "0.05254";   E = "  -0.17%  " }
    40 = @{ D = "2.918";     E = "  +0.45%  " }
    41 = @{ D = "0.5330";    E = "  +0.40%  " }
    42 = @{ D = "7.179";     E = "  -1.58%  " }
    43 = @{ D = "2.349";     E = "  -1.69%  " }
    44 = @{ D = "0.1659";    E = "  -1.91%  " }
    45 = @{ D = "8.420";     E = "  -2.48%  " }
    46 = @{ D = "0.4935";    E = "  -2.41%  " }
    47 = @{ D = "10.37";     E = "  -1.58%  " }
    48 = @{ E = "  +0.19%  " }
    49 = @{ E = "  -0.36%  " }
    50 = @{ D = "102.77";    E = "  -2.77%  " }
    51 = @{ D = "0.06266";   E = "  -1.10%  " }
}

# Rows whose new Price text parses as a plain single number (e.g. "310.27",
# "1.002", "0.000008699") must be pinned to Text format before the write,
# otherwise the interop layer will coerce the string into a floating point
# value (losing trailing zeros / switching to scientific notation). Prices
# that already contain thousands separators (e.g. "26.851.62") are never
# parsed as numbers, so they can be written directly.
$plainNumberPattern = '^-?\d+(\.\d+)?$'

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($vals.ContainsKey("D")) {
        $newPrice = $vals["D"]
        $dCell = $ws.Range("D$row")
        if ($newPrice -match $plainNumberPattern) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $newPrice
    }

    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
